$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("김혜민")
$ws.Range("A37").Value = 43798
$ws.Range("A37").Select()
